# Add new code sample: record the employee's alias change.
# Augusta Delono's alias was "AUD"; it is being changed to "ADO".
# The prior alias value is kept as a note in cell J4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J4").Value = "AUD"
$ws.Range("D4").Value = "ADO"

$ws.Range("H17").Select() | Out-Null
